$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L (year 2022) mirroring the existing K column (year 2021)
# formatting, row by row. Copy each K cell's style into the matching L
# cell first, then set the new value (and, for L8, add a number format
# that introduces a brand-new cell style).

$ws.Range("K2").Copy($ws.Range("L2"))

$ws.Range("K3").Copy($ws.Range("L3"))
$ws.Range("L3").Value = 2022

$ws.Range("K4").Copy($ws.Range("L4"))
$ws.Range("L4").Value = 370

$ws.Range("K5").Copy($ws.Range("L5"))
$ws.Range("L5").Value = 137

$ws.Range("K6").Copy($ws.Range("L6"))
$ws.Range("L6").Value = 314

$ws.Range("K7").Copy($ws.Range("L7"))
$ws.Range("L7").Value = 121

$ws.Range("K8").Copy($ws.Range("L8"))
$ws.Range("L8").Value = 50
$ws.Range("L8").NumberFormat = "#,##0"

$ws.Range("K9").Copy($ws.Range("L9"))
$ws.Range("L9").Value = 16

# Move the active selection to the top of the new column, matching the
# saved view state.
$null = $ws.Range("L2").Select()
